$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 8 (the "QNT" row) - remaining rows shift up
$ws.Rows.Item(8).Delete() | Out-Null

# Update Kaufpreis value for SOL (row 4)
$ws.Range("C4").Value = 30

# Update selection to match target workbook state
$ws.Range("C4").Select() | Out-Null
